$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1593
$ws1.Range("F6").Value  = 211
$ws1.Range("F7").Value  = 785
$ws1.Range("F8").Value  = 1065
$ws1.Range("F14").Value = 6612
$ws1.Range("F15").Value = 136
$ws1.Range("F21").Value = 15836
$ws1.Range("F22").Value = 1566
$ws1.Range("F23").Value = 26
$ws1.Range("F27").Value = 11198
$ws1.Range("F28").Value = 807
$ws1.Range("F30").Value = 279

# Sheet "全部类型" (all types) - fourth sheet, same events shifted by combined rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1593
$ws4.Range("F6").Value  = 211
$ws4.Range("F7").Value  = 785
$ws4.Range("F9").Value  = 1065
$ws4.Range("F16").Value = 6612
$ws4.Range("F17").Value = 136
$ws4.Range("F24").Value = 15836
$ws4.Range("F25").Value = 1566
$ws4.Range("F26").Value = 26
$ws4.Range("F31").Value = 11198
$ws4.Range("F32").Value = 807
$ws4.Range("F34").Value = 279
